# Update "想去人数" (attendance count) values in column F
# for worksheets "展览" and "全部类型".
# Changes:
#   F5: 848  -> 850
#   F6: 14   -> 16
#   F7: 303  -> 302
#   F8: 8593 -> 8639
#   F9: 75   -> 76
#   F19: 726 -> 729

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F5"  = 850
    "F6"  = 16
    "F7"  = 302
    "F8"  = 8639
    "F9"  = 76
    "F19" = 729
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
